# Apply the "Deploying to gh-pages" metadata refresh to the
# StructureDefinition-employee-division workbook:
#   - Metadata sheet: bump Version, Date; replace the duplicated
#     "Contact" row with real "Publisher"/"Jurisdiction" values and
#     drop the now-redundant duplicate row.
#   - Elements sheet: give the root Extension row a Division-specific
#     Short/Definition instead of the generic placeholder text.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refresh the publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$meta.Range("B9").Value = "Alvearie Team"

# The old row 10 ("Contact" / "No display for ContactDetail") becomes
# the Jurisdiction row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was an exact duplicate of the old row 10 ("Contact" / "No
# display for ContactDetail") - remove it, shifting everything below
# up by one so the sheet ends up with 20 rows instead of 21.
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# Root Extension element (row 2): Short/Definition are now specific to
# this profile instead of the generic Extension boilerplate.
$elements.Range("K2").Value = "Employee Division"
$elements.Range("L2").Value = "Code for the division of the employee"
